# Creación de funciones PWM Y MAP
# Adds LED / GPIO labels (LED ROJA, LED VERDE, LED AMARILLA and their
# corresponding GPIO pins) to the pinout sheet, highlighted with a green
# fill, and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Green fill color used to highlight the new labels (RGB 92D050 -> BGR int for COM)
$greenColor = 5296274

# Cells are populated in this exact order so that the underlying shared
# string table is built with the same indices/ordering Excel itself would
# produce for this edit.
$order = @(
    @("J13", "LED ROJA"),
    @("B18", "LED VERDE"),
    @("B19", "LED AMARILLA "),
    @("C17", "GPIO4"),
    @("C18", "GPIO2"),
    @("C19", "GPIO15"),
    @("J14", "LED VERDE "),
    @("J15", "LED AMARILLA"),
    @("B17", "LED ROJA")
)

foreach ($pair in $order) {
    $addr = $pair[0]
    $val = $pair[1]
    $r = $ws.Range($addr)
    $r.Font.Bold = $true
    $r.Interior.Color = $greenColor
    $r.Value = $val
}

# Match the new selection left behind by the edit.
$ws.Range("J13:J15").Select()
